# CTMS: Added eSign CRA Submission step
# Re-apply the SiteManagementTestData.xlsx edit via Excel COM interop.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing + new text values, in the order the new unique shared
# strings were first introduced by the original authoring session -------
$ws.Range("L5").Value  = "CRA Submission"
$ws.Range("B5").Value  = "anuragk"
$ws.Range("C5").Value  = "Password1"
$ws.Range("I4").Value  = "ReportStatusBeforeUpdate"
$ws.Range("L4").Value  = "ReportStatusAfterUpdate"
$ws.Range("L4").Font.Bold = $true
$ws.Range("F1").Value  = "SuccessMsgForRecordCreation"
$ws.Range("K4").Value  = "SuccessMsgForRecordCreation"

# Bold header cells M4:P4 (no values), matching the bold style used across
# the rest of row 4.
$ws.Range("M4").Font.Bold = $true
$ws.Range("N4").Font.Bold = $true
$ws.Range("O4").Font.Bold = $true
$ws.Range("P4").Font.Bold = $true

# --- Column widths (best effort match to the authored widths; the COM
# ColumnWidth property is quantized to 1/6-character steps internally, so
# these are the closest achievable values to the authored ones) ------------
$ws.Columns.Item(1).ColumnWidth  = 28.8333333333333
$ws.Columns.Item(2).ColumnWidth  = 8.0
$ws.Columns.Item(3).ColumnWidth  = 9.0
$ws.Columns.Item(5).ColumnWidth  = 7.0
$ws.Columns.Item(6).ColumnWidth  = 26.6666666666667
$ws.Columns.Item(7).ColumnWidth  = 21.8333333333333
$ws.Columns.Item(8).ColumnWidth  = 13.6666666666667
$ws.Columns.Item(9).ColumnWidth  = 15.6666666666667
$ws.Columns.Item(10).ColumnWidth = 27.8333333333333
$ws.Columns.Item(11).ColumnWidth = 26.6666666666667
$ws.Columns.Item(12).ColumnWidth = 13.0

# --- Selection / view -------------------------------------------------------
$ws.Range("G13").Select() | Out-Null
